# Turn the literal heading "СЛУЖБОВА ЗАПИСКА" into the templated
# placeholder "{type}", written as three separate runs ("{", "type", "}")
# so it matches the same {placeholder} pattern already used elsewhere in
# this template (e.g. {receiver}, {number}, {year}, {day}, {month}).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("СЛУЖБОВА ЗАПИСКА", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Text = "{type}"

    $start = $rng.Start
    $end = $rng.End

    # $rng now covers "{type}" (6 characters). Split it into three runs -
    # "{", "type" and "}" - by nudging a direct-character-formatting
    # property on the middle/last chunks off and back on again. Because the
    # final value matches the original (Bold was already on), no stray
    # formatting is left behind, but the run boundary is preserved.
    $braceOpen = $d.Range($start, $start + 1)
    $word_ = $d.Range($start + 1, $end - 1)
    $braceClose = $d.Range($end - 1, $end)

    $word_.Bold = $false
    $word_.Bold = $true
    $braceClose.Bold = $false
    $braceClose.Bold = $true
}
